$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.001.05"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.511.44"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Formula = '=TEXT(533.41,"0.00")'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Formula = '=TEXT(135.82,"0.00")'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "2.956.62"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "58.903.14"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").Formula = '=TEXT(22.82,"0.00")'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "2.508.69"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Formula = '=TEXT(11.04,"0.00")'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Formula = '=TEXT(323.63,"0.00")'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").Formula = '=TEXT(65.15,"0.00")'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").Formula = '=TEXT(0.420,"0.000")'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").Formula = '=TEXT(6.53,"0.00")'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = "  -2.95%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").Formula = '=TEXT(169.40,"0.00")'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  -4.01%  "

$ws.Range("D34").Formula = '=TEXT(1.38,"0.00")'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("D35").Formula = '=TEXT(18.38,"0.00")'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E35").Value = "  -1.09%  "

$ws.Range("E36").Value = "  -1.85%  "

$ws.Range("E37").Value = "  -2.91%  "

$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("E39").Value = "  -3.74%  "

$ws.Range("D40").Formula = '=TEXT(282.48,"0.00")'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("E42").Value = "  -5.16%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Formula = '=TEXT(129.64,"0.00")'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Formula = '=TEXT(10.93,"0.00")'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").Value = "1.760.97"
$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("E51").Value = "  -0.44%  "
